# Apply cryptos.xlsx price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.451.60"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.634.58"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").Value = "2.633.46"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.114.97"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "67.329.70"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "2.637.99"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "365.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +4.82%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "2.765.44"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "582.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.49%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.21%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "0.0₆0287"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.90%  "
